$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple cell updates (Price / Volume columns) ---
$ws.Range("D2").Value = '35.048.77'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").Value = '1.816.72'
$ws.Range("E3").Value = '  -0.68%  '
$ws.Range("E4").Value = '  +0.43%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.01'
$ws.Range("E5").Value = '  +2.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.614'
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.35'
$ws.Range("E8").Value = '  -8.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.330'
$ws.Range("E9").Value = '  +10.23%  '
$ws.Range("E10").Value = '  -0.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.100'
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").Value = '2.079.46'
$ws.Range("E12").Value = '  -0.58%  '
$ws.Range("D13").Value = '1.826.34'
$ws.Range("E13").Value = '  -0.02%  '
$ws.Range("E14").Value = '  -1.01%  '
$ws.Range("D17").Value = '34.977.10'
$ws.Range("E17").Value = '  +0.93%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.52'
$ws.Range("E18").Value = '  +1.72%  '
$ws.Range("D19").Value = '0.0₃0789'
$ws.Range("E19").Value = '  +0.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '239.02'
$ws.Range("E20").Value = '  -1.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.83'
$ws.Range("E21").Value = '  -2.52%  '
$ws.Range("E22").Value = '  -1.97%  '
$ws.Range("E23").Value = '  +0.40%  '
$ws.Range("E24").Value = '  +2.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '172.66'
$ws.Range("E25").Value = '  +0.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.86'
$ws.Range("E26").Value = '  -1.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.48'
$ws.Range("E27").Value = '  -1.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.121'
$ws.Range("E28").Value = '  -1.20%  '
$ws.Range("E29").Value = '  +29.21%  '
$ws.Range("E30").Value = '  +0.47%  '
$ws.Range("D31").Value = '3.337.99'
$ws.Range("E31").Value = '  +37.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0555'
$ws.Range("E32").Value = '  +5.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.96'
$ws.Range("E33").Value = '  +1.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.96'
$ws.Range("E34").Value = '  -0.21%  '
$ws.Range("E35").Value = '  -3.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.683'
$ws.Range("E38").Value = '  +2.56%  '
$ws.Range("E39").Value = '  +0.73%  '
$ws.Range("E40").Value = '  +3.56%  '
$ws.Range("D41").Value = '1.304.94'
$ws.Range("E41").Value = '  -2.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.985'
$ws.Range("E42").Value = '  +1.21%  '
$ws.Range("E43").Value = '  -4.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.60'
$ws.Range("E44").Value = '  -4.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.77'
$ws.Range("E46").Value = '  -1.71%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.36'
$ws.Range("E47").Value = '  +6.46%  '
$ws.Range("E48").Value = '  -1.40%  '
$ws.Range("D49").Value = '1.994.47'
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("E51").Value = '  +5.03%  '

# --- Row 15/16 swap: Polkadot <-> Polygon ---
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.662'
$ws.Range("E15").Value = '  +1.78%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.65'
$ws.Range("E16").Value = '  -0.56%  '

# --- Row 36/37 swap: Aave <-> TrustWalletToken ---
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.14'
$ws.Range("E36").Value = '  +6.91%  '
$ws.Range("B37").Value = 'Aave'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '92.74'
$ws.Range("E37").Value = '  +2.96%  '
